$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: add a new "UPDATE:" paragraph inside the table cell, right after
# the paragraph that ends with: for "row"). The function will show that the
# point is invalid if it is outside of this range.
# ---------------------------------------------------------------------------
$targetText = 'for "row"). The function will show that the point is invalid if it is outside of this range.'

$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*$targetText*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the paragraph ending in the point-range description."
}

$srcPara = $d.Paragraphs($targetIndex)
$srcPara.Range.InsertParagraphAfter()

# The freshly inserted (empty) paragraph is now the next one in the document.
$newPara = $d.Paragraphs($targetIndex + 1)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr>' +
  '<w:r><w:t xml:space="preserve">UPDATE: </w:t></w:r>' +
  '<w:r><w:t>The function now ensures that the destination entered is a valid building. To prevent further issues in the program there''s an additional check that prevents user from entering non-edge buildings.</w:t></w:r>' +
  '</w:p>'

$newPara.Range.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# Change 2: mark the run containing "function returns 0 to indicate
# successful validation" with a <w:lastRenderedPageBreak/> element.
# ---------------------------------------------------------------------------
$descMarker = 'function returns 0 to indicate successful validation'

$count2 = $d.Paragraphs.Count
$descIndex = -1
for ($i = 1; $i -le $count2; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*$descMarker*") {
        $descIndex = $i
        break
    }
}

if ($descIndex -eq -1) {
    throw "Could not locate the Description paragraph."
}

$descPara = $d.Paragraphs($descIndex)

$descParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="67CC5781" w14:textId="065812DF" w:rsidR="00BA6788" w:rsidRPr="000116B5" w:rsidRDefault="00BA6788" w:rsidP="00BA6788">' +
  '<w:r w:rsidRPr="00BA6788"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Description:</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r w:rsidR="00DC6D79"><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r w:rsidR="00DC6D79" w:rsidRPr="00DC6D79"><w:t>The validate function checks if an object''s weight, volume, and location meet specific standards. It gives a number as a result to show if the validation was successful or not. If the weight is not between 1 and 1000, the function returns an error code. Similarly, if the volume is not 0.25, 0.5, or 1.0, it gives an error code. Additionally, if the object''s position is outside the allowed range of coordinates, the function returns an error code. If all the parameters meet the requirements, t</w:t></w:r>' +
  '<w:r w:rsidR="000116B5"><w:t xml:space="preserve">he </w:t></w:r>' +
  '<w:r w:rsidR="00DC6D79" w:rsidRPr="00DC6D79"><w:lastRenderedPageBreak/><w:t>function returns 0 to indicate successful validation</w:t></w:r>' +
  '<w:r w:rsidR="000116B5"><w:t>.</w:t></w:r>' +
  '<w:r w:rsidR="00DC6D79"><w:br/></w:r>' +
  '</w:p>'

$descPara.Range.InsertXML($descParaXml)

Write-Host "Edit complete."
